$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-11 change from 45175 to 45183 (2023-09-06 -> 2023-09-14)
$ws.Range("C2:C11").Value = 45183
